$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people interested) counts increased slightly
# for a handful of events. The same rows appear on both the "展览" sheet
# and the "全部类型" sheet, so update F5/F8/F11/F13/F15 on both.

$updates = @{
    "F5"  = 4636
    "F8"  = 1378
    "F11" = 1079
    "F13" = 587
    "F15" = 16
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
